# Update the "Colour Code" gradient column (column B) on the "Date Colours"
# sheet. The gradient was recomputed/shifted by one row: the sequence that
# used to run from row 2 through row 41 (ending at #118dff) now runs from
# row 2 through row 44, so every existing value shifts to a new shade and
# three more rows (42-44) gain a value for the first time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Date Colours")

$colours = @{
    2  = "#fbfcff"
    3  = "#f8f9ff"
    4  = "#f4f6ff"
    5  = "#f0f4ff"
    6  = "#edf1ff"
    7  = "#e9eeff"
    8  = "#e5ebff"
    9  = "#e1e9ff"
    10 = "#dee6ff"
    11 = "#dae3ff"
    12 = "#d6e0ff"
    13 = "#d2ddff"
    14 = "#cedbff"
    15 = "#cad8ff"
    16 = "#c6d5ff"
    17 = "#c2d3ff"
    18 = "#bed0ff"
    19 = "#bacdff"
    20 = "#b6caff"
    21 = "#b2c8ff"
    22 = "#adc5ff"
    23 = "#a9c2ff"
    24 = "#a5c0ff"
    25 = "#a0bdff"
    26 = "#9cbaff"
    27 = "#97b8ff"
    28 = "#92b5ff"
    29 = "#8eb3ff"
    30 = "#89b0ff"
    31 = "#83adff"
    32 = "#7eabff"
    33 = "#79a8ff"
    34 = "#73a6ff"
    35 = "#6da3ff"
    36 = "#67a1ff"
    37 = "#619eff"
    38 = "#5a9cff"
    39 = "#5299ff"
    40 = "#4a97ff"
    41 = "#4094ff"
    42 = "#3592ff"
    43 = "#278fff"
    44 = "#118dff"
}

foreach ($row in 2..44) {
    $ws.Cells.Item($row, 2).Value = $colours[$row]
}
